# Append a new data row (row 48) to the sheet, extending the table
# from A1:D47 to A1:D48, mirroring the existing "date / weekday /
# hour / ranking" log rows already present in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds dates formatted as plain text (e.g. "2025/10/01"),
# not real Excel date serials. A leading apostrophe forces the
# COM layer to store the value as text instead of auto-converting
# it to a date serial; ClearFormats() then strips the "quote
# prefix" formatting flag that the apostrophe trick leaves behind,
# so the cell ends up as a plain, unstyled text cell - just like
# the rest of the column.
$ws.Range("A48").Value2 = "'2025/10/02"
$ws.Range("A48").ClearFormats()

# Column B holds the weekday (plain text, no special handling needed).
$ws.Range("B48").Value2 = "木"

# Columns C and D hold plain numeric values.
$ws.Range("C48").Value2 = 1
$ws.Range("D48").Value2 = 201
